$wb = $excel.ActiveWorkbook

# --- SFF sheet: fill in Snatch / Clean and Jerk results ---
$sff = $wb.Worksheets.Item("SFF")
$sff.Range("C2").Value = 45
$sff.Range("D2").Value = 67.5
$sff.Range("C3").Value = 52.5
$sff.Range("D3").Value = 72.5
$sff.Range("C4").Value = 62.5
$sff.Range("D4").Value = 77.5
$sff.Range("C5").Value = 40
$sff.Range("D5").Value = 60
$sff.Range("C6").Value = 65
$sff.Range("D6").Value = 72.5
$sff.Range("C7").Value = 70
$sff.Range("D7").Value = 90

# selection ends on D5 for SFF (active/tab-selected sheet)
$sff.Activate()
$sff.Range("D5").Select()

# --- SFM sheet: selection only moves to F6 ---
$sfm = $wb.Worksheets.Item("SFM")
$sfm.Activate()
$sfm.Range("F6").Select()

# Restore the SFF sheet as the active tab (matches tabSelected on SFF)
$sff.Activate()
